$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.432.85"
$ws.Range("E2").Value = "  -1.36%  "

$ws.Range("D3").Value = "2.637.84"
$ws.Range("E3").Value = "  -2.76%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.29%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.38"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.03%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("E8").Value = "  -3.48%  "

$ws.Range("D9").Value = "2.637.04"
$ws.Range("E9").Value = "  -2.77%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.144"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.33%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.363"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.11%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.27"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.38%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.92"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.65%  "

$ws.Range("D15").Value = "3.114.21"
$ws.Range("E15").Value = "  -2.89%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000183"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.69%  "

$ws.Range("D17").Value = "67.302.61"
$ws.Range("E17").Value = "  -1.38%  "

$ws.Range("D18").Value = "2.604.69"
$ws.Range("E18").Value = "  -4.00%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.26"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.88%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.08"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.99%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "361.16"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.37%  "

$ws.Range("E22").Value = "  -2.57%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.71"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.20%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.88"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +8.63%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.97"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.11%  "

$ws.Range("E26").Value = "  -0.01%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "70.40"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.60%  "

$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000102"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.31%  "

$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.995"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.42%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "556.93"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.97%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.95"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.97%  "

$ws.Range("E33").Value = "  -2.95%  "

$ws.Range("E35").Value = "  +3.97%  "

$ws.Range("E36").Value = "  -0.01%  "

$ws.Range("E37").Value = "  -4.69%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "157.48"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.89%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.21"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.33%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.368"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.65%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.23"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.94%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.81"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.05%  "

$ws.Range("E43").Value = "  -0.19%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.49"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.47%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.17"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.43%  "

$ws.Range("D47").Value = "0.0₆0300"
$ws.Range("E47").Value = "  -3.20%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.589"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.17%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "152.53"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.61%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.83"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.80%  "

$ws.Range("E51").Value = "  -1.38%  "
